$d = $word.ActiveDocument

# Locate the paragraph that ends the math-game notes ("Give 1 point for
# every right answer...") -- the two new notes paragraphs plus a blank
# paragraph get inserted right after it, before "Randomize several..."
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Give 1 point for every right answer*") {
        $anchor = $p
    }
}

$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p1 = $anchor.Next()

$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$p2 = $p1.Next()

$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p3 = $p2.Next()

# First new note
$p1.Range.Text = "Making a new score variable didn’t work"
$p1.Format.LeftIndent = 36

# Second new note -- the _GoBack bookmark (previously sitting after the
# "Questions regarding..." paragraph) now lands in the middle of this
# paragraph's text, splitting it into two runs.
$fullText2 = "Making calling a score counting function is too complicated. Need to find simpler way"
$p2.Range.Text = $fullText2
$p2.Format.LeftIndent = 36

$splitOffset = "Making calling a score counting function is too co".Length
$bmPos = $p2.Range.Start + $splitOffset
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# $p3 stays a blank paragraph (no indent, no text) -- matches the diff's
# lone "<w:p/>" between the new notes and "Randomize several...".
